$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1132.5454
$ws.Range("I19").Value = 1077.3334
$ws.Range("K19").Value = 1077.3334
$ws.Range("M19").Value = -902.3334
$ws.Range("H64").Value = 4465.4165
$ws.Range("I64").Value = 3998.889
$ws.Range("J64").Value = 5865
$ws.Range("K64").Value = 3998.889
$ws.Range("L64").Value = 5865
$ws.Range("M64").Value = -3750.889
$ws.Range("N64").Value = -6361
$ws.Range("H67").Value = 4465.4165
$ws.Range("I67").Value = 3998.889
$ws.Range("J67").Value = 5865
$ws.Range("K67").Value = 3998.889
$ws.Range("L67").Value = 5865
$ws.Range("M67").Value = -3140.889
$ws.Range("N67").Value = -7581
$ws.Range("H74").Value = 8473.869000000001
$ws.Range("I74").Value = 8631.772000000001
$ws.Range("K74").Value = 8631.772000000001
$ws.Range("M74").Value = -7695.772000000001
$ws.Range("H77").Value = 8473.869000000001
$ws.Range("I77").Value = 8631.772000000001
$ws.Range("K77").Value = 43158.86
$ws.Range("M77").Value = -38478.86
$ws.Range("H107").Value = 3515.7036
$ws.Range("I107").Value = 2727.8462
$ws.Range("J107").Value = 24000
$ws.Range("K107").Value = 2727.8462
$ws.Range("L107").Value = 24000
$ws.Range("M107").Value = -807.8462
$ws.Range("N107").Value = -27840
$ws.Range("H132").Value = 12587.588
$ws.Range("I132").Value = 14192.6
$ws.Range("K132").Value = 42577.8
$ws.Range("M132").Value = -40047.8
$ws.Range("H138").Value = 5054.7896
$ws.Range("I138").Value = 6541.1665
$ws.Range("J138").Value = 4368.769
$ws.Range("K138").Value = 19623.4995
$ws.Range("L138").Value = 13106.307
$ws.Range("M138").Value = -14483.4995
$ws.Range("N138").Value = -23386.307

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 163273.17
$ws.Range("I32").Value = 172532.14
$ws.Range("K32").Value = 172532.14
$ws.Range("M32").Value = -172245.14
$ws.Range("H39").Value = 5015.5
$ws.Range("I39").Value = 5015.5
$ws.Range("K39").Value = 5015.5
$ws.Range("M39").Value = -4495.5
$ws.Range("H61").Value = 3085.182
$ws.Range("I61").Value = 3093.7
$ws.Range("K61").Value = 3093.7
$ws.Range("M61").Value = -2881.7
$ws.Range("H110").Value = 2246.5
$ws.Range("I110").Value = 2862.75
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 2862.75
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -817.75
$ws.Range("N110").Value = -6090
$ws.Range("H136").Value = 3085.182
$ws.Range("I136").Value = 3093.7
$ws.Range("K136").Value = 9281.099999999999
$ws.Range("M136").Value = -6731.099999999999
$ws.Range("H138").Value = 86710
$ws.Range("J138").Value = 86710
$ws.Range("L138").Value = 86710
$ws.Range("N138").Value = -96990

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3000000
$ws.Range("I7").Value = 3000000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3000000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2999887
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 212.14285
$ws.Range("I22").Value = 198.83333
$ws.Range("K22").Value = 198.83333
$ws.Range("M22").Value = -25.83332999999999
$ws.Range("H132").Value = 98994.28999999999
$ws.Range("J132").Value = 98994.28999999999
$ws.Range("L132").Value = 98994.28999999999
$ws.Range("N132").Value = -109114.29
$ws.Range("H133").Value = 105000
$ws.Range("J133").Value = 105000
$ws.Range("L133").Value = 105000
$ws.Range("N133").Value = -115120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125400.664
$ws.Range("I16").Value = 20934.334
$ws.Range("K16").Value = 20934.334
$ws.Range("M16").Value = -20647.334
$ws.Range("H18").Value = 21999
$ws.Range("J18").Value = 21999
$ws.Range("L18").Value = 21999
$ws.Range("N18").Value = -22459
$ws.Range("H31").Value = 2047.25
$ws.Range("I31").Value = 2090.4
$ws.Range("J31").Value = 1400
$ws.Range("K31").Value = 2090.4
$ws.Range("L31").Value = 1400
$ws.Range("M31").Value = -1795.4
$ws.Range("N31").Value = -1990
$ws.Range("H34").Value = 2047.25
$ws.Range("I34").Value = 2090.4
$ws.Range("J34").Value = 1400
$ws.Range("K34").Value = 2090.4
$ws.Range("L34").Value = 1400
$ws.Range("M34").Value = -1888.4
$ws.Range("N34").Value = -1804
$ws.Range("H58").Value = 2949.681
$ws.Range("I58").Value = 2520.2334
$ws.Range("J58").Value = 3707.5293
$ws.Range("K58").Value = 2520.2334
$ws.Range("L58").Value = 3707.5293
$ws.Range("M58").Value = -2317.2334
$ws.Range("N58").Value = -4113.5293
$ws.Range("H62").Value = 4056.8572
$ws.Range("I62").Value = 3679.8
$ws.Range("K62").Value = 3679.8
$ws.Range("M62").Value = -3055.8
$ws.Range("H65").Value = 4056.8572
$ws.Range("I65").Value = 3679.8
$ws.Range("K65").Value = 18399
$ws.Range("M65").Value = -15279
$ws.Range("H68").Value = 41782.832
$ws.Range("J68").Value = 41782.832
$ws.Range("L68").Value = 41782.832
$ws.Range("N68").Value = -43280.832
$ws.Range("H71").Value = 41782.832
$ws.Range("J71").Value = 41782.832
$ws.Range("L71").Value = 125348.496
$ws.Range("N71").Value = -132836.496
$ws.Range("H113").Value = 125400.664
$ws.Range("I113").Value = 20934.334
$ws.Range("K113").Value = 20934.334
$ws.Range("M113").Value = -18764.334
$ws.Range("H114").Value = 29500
$ws.Range("J114").Value = 29500
$ws.Range("L114").Value = 29500
$ws.Range("N114").Value = -38178
$ws.Range("H132").Value = 22572.25
$ws.Range("I132").Value = 28763
$ws.Range("K132").Value = 86289
$ws.Range("M132").Value = -83759
$ws.Range("H134").Value = 3344.9443
$ws.Range("I134").Value = 3031.7
$ws.Range("K134").Value = 9095.099999999999
$ws.Range("M134").Value = -6560.099999999999
$ws.Range("H136").Value = 2949.681
$ws.Range("I136").Value = 2520.2334
$ws.Range("J136").Value = 3707.5293
$ws.Range("K136").Value = 7560.7002
$ws.Range("L136").Value = 11122.5879
$ws.Range("M136").Value = -5010.7002
$ws.Range("N136").Value = -16222.5879

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2178.2856
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2178.2856
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2178.2856
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6518.2856
$ws.Range("H126").Value = 2584.7334
$ws.Range("J126").Value = 2665
$ws.Range("L126").Value = 7995
$ws.Range("N126").Value = -12935

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2029.6
$ws.Range("I7").Value = 1288.25
$ws.Range("K7").Value = 1288.25
$ws.Range("M7").Value = -1176.25
$ws.Range("H93").Value = 6923.857
$ws.Range("J93").Value = 7618.25
$ws.Range("L93").Value = 7618.25
$ws.Range("N93").Value = -10114.25
$ws.Range("H126").Value = 2029.6
$ws.Range("I126").Value = 1288.25
$ws.Range("K126").Value = 3864.75
$ws.Range("M126").Value = -1394.75
$ws.Range("H132").Value = 3886.5789
$ws.Range("I132").Value = 3321.375
$ws.Range("K132").Value = 9964.125
$ws.Range("M132").Value = -7434.125
$ws.Range("H136").Value = 8083.722
$ws.Range("I136").Value = 4035.9167
$ws.Range("J136").Value = 16179.333
$ws.Range("K136").Value = 12107.7501
$ws.Range("L136").Value = 48537.999
$ws.Range("M136").Value = -9557.750100000001
$ws.Range("N136").Value = -53637.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20666.334
$ws.Range("I2").Value = 20666.334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 20666.334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -20554.334
$ws.Range("N2").ClearContents()
$ws.Range("H54").Value = 91250
$ws.Range("I54").Value = 101791.664
$ws.Range("K54").Value = 101791.664
$ws.Range("M54").Value = -101271.664
$ws.Range("H107").Value = 1663.5714
$ws.Range("I107").Value = 1395.3077
$ws.Range("J107").Value = 1896.0667
$ws.Range("K107").Value = 4185.9231
$ws.Range("L107").Value = 5688.2001
$ws.Range("M107").Value = -2265.9231
$ws.Range("N107").Value = -9528.2001
$ws.Range("H122").Value = 52102
$ws.Range("I122").Value = 3032.3125
$ws.Range("J122").Value = 164261.28
$ws.Range("K122").Value = 9096.9375
$ws.Range("L122").Value = 492783.84
$ws.Range("M122").Value = -6646.9375
$ws.Range("N122").Value = -497683.84
$ws.Range("H140").Value = 64685.8
$ws.Range("J140").Value = 64685.8
$ws.Range("L140").Value = 64685.8
$ws.Range("N140").Value = -75045.8
$ws.Range("H141").Value = 82607.5
$ws.Range("J141").Value = 82607.5
$ws.Range("L141").Value = 82607.5
$ws.Range("N141").Value = -92967.5
